# Apply the edits described by the diff to the 乳制品 (dairy products) sheet:
#   1. For each year block (rows grouped in 4s: quarters A,B,C,D), swap the
#      content of the "B" quarter row and the "C" quarter row (columns A:E).
#   2. Delete columns F and G entirely (they are dropped from the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data starts at row 2 and runs through row 81, in blocks of 4 rows per year
# (A, B, C, D quarters). Swap the 2nd (B) and 3rd (C) row of every block.
for ($blockStart = 2; $blockStart -le 81; $blockStart += 4) {
    $bRow = $blockStart + 1
    $cRow = $blockStart + 2

    $bVals = $ws.Range("A$bRow`:E$bRow").Value2
    $cVals = $ws.Range("A$cRow`:E$cRow").Value2

    $ws.Range("A$bRow`:E$bRow").Value2 = $cVals
    $ws.Range("A$cRow`:E$cRow").Value2 = $bVals
}

# Remove the now-unwanted columns F (乳制品产销率) and G (乳制品销售量).
$ws.Range("F1:G81").EntireColumn.Delete()
